$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "P1-2B": move the -180 entry from B12 (Cash) to G12 (Accounts Payable)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("P1-2B")
$ws2.Range("B12").Value = ""
$ws2.Range("G12").Value = 180
$ws2.Range("A12").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "P1-3B": re-enter the transaction data for rows 4-12
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("P1-3B")

# Row 4 - clear out all the old entries
$ws3.Range("B4").Value = ""
$ws3.Range("C4").Value = ""
$ws3.Range("D4").Value = ""
$ws3.Range("E4").Value = ""
$ws3.Range("G4").Value = ""
$ws3.Range("H4").Value = ""
$ws3.Range("I4").Value = ""

# Row 5 - new values
$ws3.Range("B5").Value = 8000
$ws3.Range("C5").Value = 4000
$ws3.Range("D5").Value = 1300
$ws3.Range("E5").Value = 25000
$ws3.Range("G5").Value = 14400
$ws3.Range("I5").Value = 5300
$ws3.Range("J5").Value = -2600

# Row 6 - cleared
$ws3.Range("B6").Value = ""
$ws3.Range("G6").Value = ""

# Row 7 - cleared
$ws3.Range("B7").Value = ""
$ws3.Range("C7").Value = ""
$ws3.Range("I7").Value = ""

# Row 8 - cleared
$ws3.Range("B8").Value = ""
$ws3.Range("E8").Value = ""
$ws3.Range("G8").Value = ""
$ws3.Range("J8").Value = ""

# Row 9 - cleared
$ws3.Range("B9").Value = ""
$ws3.Range("K9").Value = ""

# Row 10 - cleared
$ws3.Range("B10").Value = ""

# Row 11 - cleared
$ws3.Range("B11").Value = ""
$ws3.Range("G11").Value = ""

# Row 12 - cleared
$ws3.Range("B12").Value = ""
$ws3.Range("J12").Value = ""

$ws3.Range("B5").Select() | Out-Null
$ws3.Activate() | Out-Null
